$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 5, 6, 9 have their "variable" columns (D, H, J, K, L, M, N, O, P, Q)
# cyclically re-shuffled among each other:
#   new row2 <- old row9
#   new row5 <- old row6
#   new row6 <- old row2
#   new row9 <- old row5
# Capture the original values first (using Value() which properly invokes the
# getter in this environment), then write them back in the new order.

function Get-RowData($row) {
    return [PSCustomObject]@{
        D = $ws.Range("D$row").Value()
        H = $ws.Range("H$row").Value()
        J = $ws.Range("J$row").Value()
        K = $ws.Range("K$row").Value()
        L = $ws.Range("L$row").Value()
        M = $ws.Range("M$row").Value()
        N = $ws.Range("N$row").Value()
        O = $ws.Range("O$row").Value()
        P = $ws.Range("P$row").Value()
        Q = $ws.Range("Q$row").Value()
    }
}

function Set-RowData($row, $data) {
    $ws.Range("D$row").Value = $data.D
    $ws.Range("H$row").Value = $data.H
    $ws.Range("J$row").Value = $data.J
    $ws.Range("K$row").Value = $data.K
    $ws.Range("L$row").Value = $data.L
    $ws.Range("M$row").Value = $data.M
    $ws.Range("N$row").Value = $data.N
    $ws.Range("O$row").Value = $data.O
    $ws.Range("P$row").Value = $data.P
    $ws.Range("Q$row").Value = $data.Q
}

$row2 = Get-RowData 2
$row5 = Get-RowData 5
$row6 = Get-RowData 6
$row9 = Get-RowData 9

Set-RowData 2 $row9
Set-RowData 5 $row6
Set-RowData 6 $row2
Set-RowData 9 $row5
